$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# Remove the now-obsolete last data row in "Recommandations" (old row 32, SICOR CI)
# so the sheet shrinks from A1:G32 to A1:G31, matching the refreshed BRVM feed.
$ws1.Rows.Item(32).Delete()

# --- Sheet "Recommandations": refresh rows 2-31 with the latest BRVM figures ---
# Row 2
$ws1.Range("A2").Value = 'BRVM - CONSOMMATION DISCRETIONNAIRE'
$ws1.Range("B2").Value = 0
$ws1.Range("C2").Value = 4
$ws1.Range("D2").Value = 708.18
$ws1.Range("E2").Value = 181.8
$ws1.Range("F2").Value = '🟡 Observer'
$ws1.Range("G2").Value = '➖ Neutre'

# Row 3
$ws1.Range("A3").Value = 'BRVM - SERVICES FINANCIERS'
$ws1.Range("B3").Value = 0
$ws1.Range("C3").Value = 4
$ws1.Range("D3").Value = 607.02
$ws1.Range("E3").Value = 152.66
$ws1.Range("F3").Value = '🟡 Observer'
$ws1.Range("G3").Value = '➖ Neutre'

# Row 4
$ws1.Range("A4").Value = 'BRVM - INDUSTRIELS'
$ws1.Range("B4").Value = 0
$ws1.Range("C4").Value = 4
$ws1.Range("D4").Value = 589.07
$ws1.Range("E4").Value = 149.68
$ws1.Range("F4").Value = '🟡 Observer'
$ws1.Range("G4").Value = '➖ Neutre'

# Row 5
$ws1.Range("A5").Value = 'BRVM-PRESTIGE'
$ws1.Range("B5").Value = 0
$ws1.Range("C5").Value = 4
$ws1.Range("D5").Value = 582.6
$ws1.Range("E5").Value = 146.08
$ws1.Range("F5").Value = '🟡 Observer'
$ws1.Range("G5").Value = '➖ Neutre'

# Row 6
$ws1.Range("A6").Value = 'BRVM - SERVICES PUBLICS'
$ws1.Range("B6").Value = 0
$ws1.Range("C6").Value = 4
$ws1.Range("D6").Value = 471.47
$ws1.Range("E6").Value = 116.89
$ws1.Range("F6").Value = '🟡 Observer'
$ws1.Range("G6").Value = '➖ Neutre'

# Row 7
$ws1.Range("A7").Value = 'BRVM - ENERGIE'
$ws1.Range("B7").Value = 0
$ws1.Range("C7").Value = 4
$ws1.Range("D7").Value = 465.14
$ws1.Range("E7").Value = 115.34
$ws1.Range("F7").Value = '🟡 Observer'
$ws1.Range("G7").Value = '➖ Neutre'

# Row 8
$ws1.Range("A8").Value = 'BRVM - TELECOMMUNICATIONS'
$ws1.Range("B8").Value = 0
$ws1.Range("C8").Value = 4
$ws1.Range("D8").Value = 376.94
$ws1.Range("E8").Value = 94.29
$ws1.Range("F8").Value = '🟡 Observer'
$ws1.Range("G8").Value = '➖ Neutre'

# Row 9
$ws1.Range("A9").Value = 'TRACTAFRIC MOTORS CI (PRSC)'
$ws1.Range("B9").Value = 3
$ws1.Range("C9").Value = 1
$ws1.Range("D9").Value = 18.54
$ws1.Range("E9").Value = 6.32
$ws1.Range("F9").Value = '🟢 Achat'
$ws1.Range("G9").Value = '✅ Renforcer'

# Row 10
$ws1.Range("A10").Value = 'UNIWAX CI (UNXC)'
$ws1.Range("B10").Value = 2
$ws1.Range("C10").Value = 1
$ws1.Range("D10").Value = 9.01
$ws1.Range("E10").Value = -3.79
$ws1.Range("F10").Value = '🟡 Observer'
$ws1.Range("G10").Value = '👀 À surveiller'

# Row 11
$ws1.Range("A11").Value = 'SUCRIVOIRE (SCRC)'
$ws1.Range("B11").Value = 2
$ws1.Range("C11").Value = 1
$ws1.Range("D11").Value = 8.15
$ws1.Range("E11").Value = -6.32
$ws1.Range("F11").Value = '🟡 Observer'
$ws1.Range("G11").Value = '👀 À surveiller'

# Row 12
$ws1.Range("A12").Value = 'CFAO MOTORS CI (CFAC)'
$ws1.Range("B12").Value = 1
$ws1.Range("C12").Value = 0
$ws1.Range("D12").Value = 7.19
$ws1.Range("E12").Value = 7.19
$ws1.Range("F12").Value = '🟡 Observer'
$ws1.Range("G12").Value = '➖ Neutre'

# Row 13
$ws1.Range("A13").Value = 'SAFCA CI (SAFC)'
$ws1.Range("B13").Value = 1
$ws1.Range("C13").Value = 0
$ws1.Range("D13").Value = 5.57
$ws1.Range("E13").Value = 5.57
$ws1.Range("F13").Value = '🟡 Observer'
$ws1.Range("G13").Value = '➖ Neutre'

# Row 14
$ws1.Range("A14").Value = 'ORAGROUP TOGO (ORGT)'
$ws1.Range("B14").Value = 1
$ws1.Range("C14").Value = 0
$ws1.Range("D14").Value = 5
$ws1.Range("E14").Value = 5
$ws1.Range("F14").Value = '🟡 Observer'
$ws1.Range("G14").Value = '➖ Neutre'

# Row 15
$ws1.Range("A15").Value = 'ECOBANK TRANS. INCORP. TG (ETIT)'
$ws1.Range("B15").Value = 2
$ws1.Range("C15").Value = 1
$ws1.Range("D15").Value = 4.73
$ws1.Range("E15").Value = -4.17
$ws1.Range("F15").Value = '🟡 Observer'
$ws1.Range("G15").Value = '👀 À surveiller'

# Row 16
$ws1.Range("A16").Value = 'BANK OF AFRICA BN (BOAB)'
$ws1.Range("B16").Value = 1
$ws1.Range("C16").Value = 0
$ws1.Range("D16").Value = 3.74
$ws1.Range("E16").Value = 3.74
$ws1.Range("F16").Value = '🟡 Observer'
$ws1.Range("G16").Value = '➖ Neutre'

# Row 17
$ws1.Range("A17").Value = 'ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)'
$ws1.Range("B17").Value = 1
$ws1.Range("C17").Value = 0
$ws1.Range("D17").Value = 2.86
$ws1.Range("E17").Value = 2.86
$ws1.Range("F17").Value = '🟡 Observer'
$ws1.Range("G17").Value = '➖ Neutre'

# Row 18
$ws1.Range("A18").Value = 'AFRICA GLOBAL LOGISTICS CI (SDSC)'
$ws1.Range("B18").Value = 1
$ws1.Range("C18").Value = 0
$ws1.Range("D18").Value = 2.58
$ws1.Range("E18").Value = 2.58
$ws1.Range("F18").Value = '🟡 Observer'
$ws1.Range("G18").Value = '➖ Neutre'

# Row 19
$ws1.Range("A19").Value = 'FILTISAC CI (FTSC)'
$ws1.Range("B19").Value = 1
$ws1.Range("C19").Value = 1
$ws1.Range("D19").Value = 0.82
$ws1.Range("E19").Value = 3.68
$ws1.Range("F19").Value = '🟡 Observer'
$ws1.Range("G19").Value = '👀 À surveiller'

# Row 20
$ws1.Range("A20").Value = 'VIVO ENERGY CI (SHEC)'
$ws1.Range("B20").Value = 1
$ws1.Range("C20").Value = 2
$ws1.Range("D20").Value = 0.44
$ws1.Range("E20").Value = -1.58
$ws1.Range("F20").Value = '🟡 Observer'
$ws1.Range("G20").Value = '👀 À surveiller'

# Row 21
$ws1.Range("A21").Value = 'LOTERIE NATIONALE DU BENIN (LNBB)'
$ws1.Range("B21").Value = 1
$ws1.Range("C21").Value = 1
$ws1.Range("D21").Value = 0.06
$ws1.Range("E21").Value = 2.56
$ws1.Range("F21").Value = '🟡 Observer'
$ws1.Range("G21").Value = '👀 À surveiller'

# Row 22
$ws1.Range("A22").Value = 'SOCIETE IVOIRIENNE DE BANQUE  (SIBC)'
$ws1.Range("B22").Value = 1
$ws1.Range("C22").Value = 1
$ws1.Range("D22").Value = 0.01
$ws1.Range("E22").Value = 3.13
$ws1.Range("F22").Value = '🟡 Observer'
$ws1.Range("G22").Value = '👀 À surveiller'

# Row 23
$ws1.Range("A23").Value = 'SAPH CI (SPHC)'
$ws1.Range("B23").Value = 0
$ws1.Range("C23").Value = 1
$ws1.Range("D23").Value = -0.98
$ws1.Range("E23").Value = -0.98
$ws1.Range("F23").Value = '🟡 Observer'
$ws1.Range("G23").Value = '➖ Neutre'

# Row 24
$ws1.Range("A24").Value = 'BANK OF AFRICA NG (BOAN)'
$ws1.Range("B24").Value = 0
$ws1.Range("C24").Value = 1
$ws1.Range("D24").Value = -1.14
$ws1.Range("E24").Value = -1.14
$ws1.Range("F24").Value = '🟡 Observer'
$ws1.Range("G24").Value = '➖ Neutre'

# Row 25
$ws1.Range("A25").Value = 'ONATEL BF (ONTBF)'
$ws1.Range("B25").Value = 0
$ws1.Range("C25").Value = 1
$ws1.Range("D25").Value = -1.15
$ws1.Range("E25").Value = -1.15
$ws1.Range("F25").Value = '🟡 Observer'
$ws1.Range("G25").Value = '➖ Neutre'

# Row 26
$ws1.Range("A26").Value = 'CIE CI (CIEC)'
$ws1.Range("B26").Value = 0
$ws1.Range("C26").Value = 1
$ws1.Range("D26").Value = -1.85
$ws1.Range("E26").Value = -1.85
$ws1.Range("F26").Value = '🟡 Observer'
$ws1.Range("G26").Value = '➖ Neutre'

# Row 27
$ws1.Range("A27").Value = 'SMB CI (SMBC)'
$ws1.Range("B27").Value = 0
$ws1.Range("C27").Value = 1
$ws1.Range("D27").Value = -2.15
$ws1.Range("E27").Value = -2.15
$ws1.Range("F27").Value = '🟡 Observer'
$ws1.Range("G27").Value = '➖ Neutre'

# Row 28
$ws1.Range("A28").Value = 'NESTLE CI (NTLC)'
$ws1.Range("B28").Value = 0
$ws1.Range("C28").Value = 2
$ws1.Range("D28").Value = -3.48
$ws1.Range("E28").Value = -2.59
$ws1.Range("F28").Value = '🟡 Observer'
$ws1.Range("G28").Value = '➖ Neutre'

# Row 29
$ws1.Range("A29").Value = 'SICOR CI (SICC)'
$ws1.Range("B29").Value = 0
$ws1.Range("C29").Value = 1
$ws1.Range("D29").Value = -5.26
$ws1.Range("E29").Value = -5.26
$ws1.Range("F29").Value = '🟡 Observer'
$ws1.Range("G29").Value = '➖ Neutre'

# Row 30
$ws1.Range("A30").Value = 'BERNABE CI (BNBC)'
$ws1.Range("B30").Value = 1
$ws1.Range("C30").Value = 2
$ws1.Range("D30").Value = -5.94
$ws1.Range("E30").Value = -7.23
$ws1.Range("F30").Value = '🟡 Observer'
$ws1.Range("G30").Value = '👀 À surveiller'

# Row 31
$ws1.Range("A31").Value = 'UNILEVER CI (UNLC)'
$ws1.Range("B31").Value = 0
$ws1.Range("C31").Value = 1
$ws1.Range("D31").Value = -7.48
$ws1.Range("E31").Value = -7.48
$ws1.Range("F31").Value = '🟡 Observer'
$ws1.Range("G31").Value = '➖ Neutre'

# --- Sheet "Top_YTD": refresh rows 2-11 with the latest YTD progression figures ---
$ws2.Range("A2").Value = 'BRVM - CONSOMMATION DISCRETIONNAIRE'
$ws2.Range("B2").Value = 5787.21
$ws2.Range("A3").Value = 'BRVM - SERVICES FINANCIERS'
$ws2.Range("B3").Value = 3917.04
$ws2.Range("A4").Value = 'BRVM - INDUSTRIELS'
$ws2.Range("B4").Value = 3637.95
$ws2.Range("A5").Value = 'BRVM-PRESTIGE'
$ws2.Range("B5").Value = 3541.38
$ws2.Range("A6").Value = 'BRVM - SERVICES PUBLICS'
$ws2.Range("B6").Value = 2153.01
$ws2.Range("A7").Value = 'BRVM - ENERGIE'
$ws2.Range("B7").Value = 2088.24
$ws2.Range("A8").Value = 'BRVM - TELECOMMUNICATIONS'
$ws2.Range("B8").Value = 1323.34
$ws2.Range("A9").Value = 'TRACTAFRIC MOTORS CI (PRSC)'
$ws2.Range("B9").Value = 19.47
$ws2.Range("A10").Value = 'UNIWAX CI (UNXC)'
$ws2.Range("B10").Value = 8.91
$ws2.Range("A11").Value = 'SUCRIVOIRE (SCRC)'
$ws2.Range("B11").Value = 7.73

Write-Host "Update complete"